$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "19.977.08"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "1.413.60"
$ws.Range("E3").Value = "  -2.37%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("E5").Value = "  -0.65%  "
Set-TextValue $ws.Range("D6") "275.78"
$ws.Range("E6").Value = "  -0.76%  "
Set-TextValue $ws.Range("D7") "0.3680"
$ws.Range("E7").Value = "  -1.05%  "
Set-TextValue $ws.Range("D8") "0.3104"
$ws.Range("E8").Value = "  +0.31%  "
Set-TextValue $ws.Range("D9") "39.87"
$ws.Range("E9").Value = "  -2.74%  "
Set-TextValue $ws.Range("D10") "1.035"
$ws.Range("E10").Value = "  +2.30%  "
Set-TextValue $ws.Range("D11") "0.06491"
$ws.Range("E11").Value = "  -1.50%  "
Set-TextValue $ws.Range("D12") "1.000"
$ws.Range("E12").Value = "  -0.66%  "
Set-TextValue $ws.Range("D13") "5.467"
$ws.Range("E13").Value = "  +0.73%  "
Set-TextValue $ws.Range("D14") "17.57"
$ws.Range("E14").Value = "  +1.72%  "
Set-TextValue $ws.Range("D15") "6.171"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "1.412.79"
$ws.Range("E16").Value = "  -2.76%  "
Set-TextValue $ws.Range("D17") "0.00001016"
$ws.Range("E17").Value = "  -0.38%  "
Set-TextValue $ws.Range("D18") "0.05674"
$ws.Range("E18").Value = "  -3.63%  "
Set-TextValue $ws.Range("D19") "1.000"
$ws.Range("E19").Value = "  -0.69%  "
Set-TextValue $ws.Range("D20") "70.80"
$ws.Range("E20").Value = "  -8.29%  "
Set-TextValue $ws.Range("D21") "5.585"
$ws.Range("E21").Value = "  -2.57%  "
Set-TextValue $ws.Range("D22") "14.69"
$ws.Range("E22").Value = "  +1.13%  "
Set-TextValue $ws.Range("D23") "10.97"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "19.982.20"
$ws.Range("E25").Value = "  -2.34%  "
Set-TextValue $ws.Range("D26") "2.276"
$ws.Range("E26").Value = "  +1.76%  "
Set-TextValue $ws.Range("D27") "132.90"
$ws.Range("E27").Value = "  -7.14%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "1.571.84"
$ws.Range("E29").Value = "  -2.85%  "
Set-TextValue $ws.Range("D30") "109.92"
$ws.Range("E30").Value = "  +0.59%  "
Set-TextValue $ws.Range("D31") "3.894"
$ws.Range("E31").Value = "  +15.56%  "
Set-TextValue $ws.Range("D32") "5.185"
$ws.Range("E32").Value = "  -4.07%  "
Set-TextValue $ws.Range("D33") "0.8107"
$ws.Range("E33").Value = "  -12.04%  "
Set-TextValue $ws.Range("D34") "0.07761"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +3.91%  "
Set-TextValue $ws.Range("D36") "4.882"
$ws.Range("E36").Value = "  +2.81%  "
Set-TextValue $ws.Range("D37") "0.05810"
$ws.Range("E37").Value = "  +3.12%  "
Set-TextValue $ws.Range("D38") "8.041"
$ws.Range("E38").Value = "  -3.67%  "
Set-TextValue $ws.Range("D39") "1.000"
$ws.Range("E39").Value = "  -0.61%  "
Set-TextValue $ws.Range("D40") "0.02046"
$ws.Range("E40").Value = "  +0.07%  "
Set-TextValue $ws.Range("D41") "10.42"
$ws.Range("E41").Value = "  -5.02%  "
Set-TextValue $ws.Range("D42") "0.1879"
$ws.Range("E42").Value = "  -1.70%  "
Set-TextValue $ws.Range("D43") "1.098"
$ws.Range("E43").Value = "  -2.69%  "
Set-TextValue $ws.Range("D44") "0.5288"
$ws.Range("E44").Value = "  -0.79%  "
Set-TextValue $ws.Range("D47") "117.30"
$ws.Range("E47").Value = "  +5.75%  "
Set-TextValue $ws.Range("D48") "0.5174"
$ws.Range("E48").Value = "  +0.28%  "
Set-TextValue $ws.Range("D49") "1.764"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("E51").Value = "  -0.52%  "

# Row 45/46: coin order swapped (EnergySwap <-> PancakeSwap) with updated price/volume
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D45") "3.532"
$ws.Range("E45").Value = "  -1.74%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "12.22"
$ws.Range("E46").Value = "  +1.29%  "
